$d = $word.ActiveDocument

# 1. "... Sơn ON PROGRESS" -> "... DONE" (map selection note)
$d.Content.Find.Execute("Sơn ON PROGRESS", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "DONE", 2) | Out-Null

# 2. "Tạo thêm 2 map cho game." -> append " DONE"
$d.Content.Find.Execute("Tạo thêm 2 map cho game.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Tạo thêm 2 map cho game. DONE", 2) | Out-Null

# 3. "Thêm quân." -> append " Cường ON PROGRESS"
$d.Content.Find.Execute("Thêm quân.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Thêm quân. Cường ON PROGRESS", 2) | Out-Null

# 4. "Rearrange menu game." -> append " Hoàng ON PROGRESS"
$d.Content.Find.Execute("Rearrange menu game.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Rearrange menu game. Hoàng ON PROGRESS", 2) | Out-Null

# 5. "Resize game để fit màn hình 4 inch." paragraph: append new text run after the tab
$pResize = $d.Paragraphs(8)
$rResize = $pResize.Range
$rResize.Collapse(0)
$rResize.InsertAfter("Sơn ON PROGRESS")

# 6. Logo paragraph -> append " DONE " (trailing space preserved)
$d.Content.Find.Execute("Làm hoặc kiếm 1 cái logo cho Main Screen với chữ Canyon Defense.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Làm hoặc kiếm 1 cái logo cho Main Screen với chữ Canyon Defense. DONE ", 2) | Out-Null

# 7. Last (empty) paragraph -> becomes a bulleted list item (reuse numId=1) with text, keep spacing
$pLogo = $d.Paragraphs(10)
$pTrail = $d.Paragraphs(11)
$mark = $d.Range($pLogo.Range.End, $pTrail.Range.Start)
$mark.Delete()

$pLogo2 = $d.Paragraphs(10)
$rSplit = $pLogo2.Range
$rSplit.Collapse(0)
$rSplit.InsertParagraphAfter()

$pNew = $d.Paragraphs(11)
$pNew.SpaceAfter = 10
$pNew.SpaceBefore = 0
$pNew.Range.InsertAfter("Lưu score cho cái nút Score (optional)")

# 8. Section docGrid charSpace 8192 -> 12288
$d.Sections(1).PageSetup.CharsLine = 0
